$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password placeholder text in C2
$ws.Range("C2").Value = "Enter your password here"

# Move the active selection to C6
$ws.Range("C6").Select()

# Adjust column widths for B and C
$ws.Columns.Item(2).ColumnWidth = 21.97
$ws.Columns.Item(3).ColumnWidth = 18.93
